$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Populate F15 first so its shared string is registered before F14's,
# matching the shared string table order produced by the original edit.
$ws.Range("F15").Value = "Learning and implementing websockets and finishing message deletion"
$ws.Range("F14").Value = "Scrum, working on message deletion backend and learning about websockets"
$ws.Range("F16").Value = "Scrum, meetings, implementing backend for panic buttons"

$ws.Range("E14").Value = 7
$ws.Range("E15").Value = 8
$ws.Range("E16").Value = 7

$ws.Range("F16").Select()
